# Update cryptocurrency price/volume data (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.729.72"
$ws.Range("D3").Value = "1.904.63"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'311.95"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.5212"
$ws.Range("E7").Value = "  +5.48%  "
$ws.Range("D8").Value = "'0.3785"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'21.23"
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("D11").Value = "'0.9032"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "'0.07660"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.911.11"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "'5.449"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'92.13"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "'0.000008707"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'0.9998"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "27.766.38"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "'5.138"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "2.151.02"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'6.626"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'153.53"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'1.868"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'2.157"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "'114.65"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'4.856"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'0.09027"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'3.183"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "'4.841"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "'1.230"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "'0.7801"
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("D36").Value = "'0.02094"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").Value = "'2.599"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").Value = "'3.077"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").Value = "'1.092"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").Value = "'0.5559"
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("D41").Value = "'0.05284"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'6.725"
$ws.Range("D43").Value = "'114.89"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "'8.518"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "'10.50"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'0.9989"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'1.616"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "'66.77"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "'0.06007"
$ws.Range("E51").Value = "  -0.78%  "
